$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B/C values (rows 2-49) and add new row 50 with recomputed simulation values.
$values = @(
    @(2, 3.414560184563785, 0.9715082175781996),
    @(3, 6.873669115625596, 2.009117360810823),
    @(4, 6.963847387791335, 3.09163508699771),
    @(5, 11.36172526641021, 4.076103905791267),
    @(6, 13.87160446187429, 5.01374544178367),
    @(7, 14.28838018689751, 6.115755149103618),
    @(8, 15.03824802640578, 7.140384260857716),
    @(9, 16.17811102956089, 8.062869346176978),
    @(10, 20.24774837571299, 9.421480561871418),
    @(11, 21.7086888204447, 10.51270186825795),
    @(12, 23.4117721480372, 12.26767129716399),
    @(13, 23.88440416133724, 13.41749628068702),
    @(14, 26.98253972628149, 14.66708669598139),
    @(15, 27.38081924812407, 15.77205972740558),
    @(16, 28.29669026642373, 16.8887124757207),
    @(17, 28.7434193179763, 17.92246773666597),
    @(18, 28.9354713203076, 19.10727674593277),
    @(19, 30.6095886493334, 20.15015421034203),
    @(20, 36.73907573072835, 21.14947648840182),
    @(21, 41.69392564439558, 22.4760546165044),
    @(22, 42.93556203367105, 23.85725923425377),
    @(23, 48.32865408425954, 25.02827306920507),
    @(24, 50.1861639085264, 26.15889293224559),
    @(25, 58.54405620326673, 27.39019954217174),
    @(26, 64.09129550975992, 28.42418790707073),
    @(27, 66.14449522269737, 29.62294798307177),
    @(28, 68.73700433170987, 30.59390238201974),
    @(29, 68.8606017629418, 31.89321862087498),
    @(30, 71.50574123286495, 32.83012823301632),
    @(31, 75.12712002498056, 33.81195371436943),
    @(32, 76.28568609339311, 34.82135019108392),
    @(33, 77.38559530561498, 36.15849802342053),
    @(34, 79.55035189048559, 37.33434630705788),
    @(35, 79.77715925669702, 38.45041929202027),
    @(36, 80.4089719186379, 39.54216948758582),
    @(37, 82.62247651712045, 40.58430055430698),
    @(38, 83.05855086897253, 41.6678574096573),
    @(39, 84.06118153322561, 42.84465069527783),
    @(40, 84.76166640437168, 43.85226849392561),
    @(41, 89.92301760418381, 44.90498511600848),
    @(42, 90.34337152317578, 46.06698466918774),
    @(43, 91.23291768875649, 47.13258111041165),
    @(44, 91.34110983542263, 48.31383260911085),
    @(45, 92.74486968610623, 49.57902081297773),
    @(46, 93.30768170106533, 50.6718462794841),
    @(47, 96.04071387289376, 51.73591260438401),
    @(48, 99.00351516661436, 52.89390434487999),
    @(49, 99.08607720026427, 53.98076850960503),
    @(50, 99.63213412225934, 55.15369795941399)
)

foreach ($row in $values) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $ws.Cells.Item($r, 2).Value2 = $b
    $ws.Cells.Item($r, 3).Value2 = $c
}

# Row 50 is new: set the A column index value (48) and copy A49 formatting so the style matches the rest of column A.
$ws.Cells.Item(50, 1).Value2 = 48
$ws.Cells.Item(49, 1).Copy()
$ws.Cells.Item(50, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
